# chemistry update to 20190730
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14
$ws.Range("A14").Value = 20190729
$ws.Range("B14").Value = "day 5"
$ws.Range("C14").Value = 0.54166666666666663
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = "T2"
$ws.Range("F14").Value = 51100

# Row 15
$ws.Range("A15").Value = 20190729
$ws.Range("B15").Value = "day 5"
$ws.Range("C15").Value = 0.54166666666666663
$ws.Range("D15").Value = 17
$ws.Range("E15").Value = "T1"
$ws.Range("F15").Value = 53600

# Row 16
$ws.Range("A16").Value = 20190729
$ws.Range("B16").Value = "day 5"
$ws.Range("C16").Value = 0.54166666666666696
$ws.Range("D16").Value = 29
$ws.Range("E16").Value = "T0"
$ws.Range("F16").Value = 56400

# Row 17
$ws.Range("A17").Value = 20190729
$ws.Range("B17").Value = "day 5"
$ws.Range("C17").Value = 0.54166666666666696
$ws.Range("D17").Value = "algae.header"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = 916000

# Row 18
$ws.Range("A18").Value = 20190730
$ws.Range("B18").Value = "day 6"
$ws.Range("C18").Value = 0.45833333333333331
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = "T2"
$ws.Range("F18").Value = 56400

# Row 19
$ws.Range("A19").Value = 20190730
$ws.Range("B19").Value = "day 6"
$ws.Range("C19").Value = 0.45833333333333331
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = "T1"
$ws.Range("F19").Value = 59200

# Row 20
$ws.Range("A20").Value = 20190730
$ws.Range("B20").Value = "day 6"
$ws.Range("C20").Value = 0.45833333333333331
$ws.Range("D20").Value = 29
$ws.Range("E20").Value = "T0"
$ws.Range("F20").Value = 35000

# Row 21
$ws.Range("A21").Value = 20190730
$ws.Range("B21").Value = "day 6"
$ws.Range("C21").Value = 0.45833333333333331
$ws.Range("D21").Value = "algae.header"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = 462000

# Apply time number format to the new C column cells (matches existing style s="1")
$ws.Range("C14:C21").NumberFormat = "h:mm"

# Update the sheet view: scrolled down a bit and selection on F19
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F19").Select()
